$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting used by the previous row (A6) onto the new row's date cell
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# Fill in the new timetable entry
$ws.Range("A7").Value = 43927
$ws.Range("B7").Value = "18:00-20:30"
$ws.Range("C7").Value = "Discussion about progrees last week, writing meeting minute"

# Restore the active selection to match the saved workbook state
$ws.Range("D10").Select()
